$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$ws.Unprotect()

$ws.Range("B3").Value = "Дата: 24.08.2021"
$ws.Range("E21").Value = 77762
$ws.Range("E23").Value = 120000
$ws.Range("E24").Value = 742238
